$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# 112 cell updates
$updates = @(
    ,@('D2', '59.123.52')
    ,@('E2', '  +0.63%  ')
    ,@('D3', '2.672.91')
    ,@('E3', '  +4.61%  ')
    ,@('E4', '  +0.24%  ')
    ,@('D5', '517.47')
    ,@('E5', '  +2.62%  ')
    ,@('D6', '145.07')
    ,@('E6', '  +2.35%  ')
    ,@('D7', '0.998')
    ,@('E7', '  -0.18%  ')
    ,@('D8', '0.568')
    ,@('E8', '  +3.05%  ')
    ,@('D9', '2.708.89')
    ,@('E9', '  +5.86%  ')
    ,@('D10', '6.23')
    ,@('E10', '  +0.60%  ')
    ,@('D11', '0.107')
    ,@('E11', '  +6.60%  ')
    ,@('D12', '0.337')
    ,@('E12', '  +2.60%  ')
    ,@('E13', '  -0.75%  ')
    ,@('D14', '3.176.56')
    ,@('E14', '  +5.68%  ')
    ,@('D15', '59.120.31')
    ,@('E15', '  +0.62%  ')
    ,@('D16', '21.08')
    ,@('E16', '  +2.98%  ')
    ,@('D17', '0.0000138')
    ,@('E17', '  +3.36%  ')
    ,@('D18', '2.690.87')
    ,@('E18', '  +5.01%  ')
    ,@('B19', 'Polkadot')
    ,@('C19', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot')
    ,@('D19', '4.55')
    ,@('E19', '  +1.41%  ')
    ,@('B20', 'BitcoinCash')
    ,@('C20', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch')
    ,@('D20', '346.16')
    ,@('E20', '  +4.76%  ')
    ,@('D21', '10.47')
    ,@('E21', '  +4.72%  ')
    ,@('D22', '6.20')
    ,@('E22', '  +5.02%  ')
    ,@('E23', '  +0.11%  ')
    ,@('D24', '61.08')
    ,@('E24', '  +2.93%  ')
    ,@('D25', '0.421')
    ,@('E25', '  +4.46%  ')
    ,@('D26', '2.783.44')
    ,@('E26', '  +3.93%  ')
    ,@('D27', '0.990')
    ,@('E27', '  -0.86%  ')
    ,@('D28', '0.161')
    ,@('E28', '  +2.39%  ')
    ,@('D29', '0.0₃0819')
    ,@('E29', '  +6.37%  ')
    ,@('D30', '7.23')
    ,@('E30', '  +6.10%  ')
    ,@('D31', '1.00')
    ,@('E31', '  -0.04%  ')
    ,@('D32', '6.43')
    ,@('E32', '  +11.46%  ')
    ,@('D33', '19.09')
    ,@('E33', '  +3.34%  ')
    ,@('D34', '1.58')
    ,@('E34', '  +3.28%  ')
    ,@('D35', '150.26')
    ,@('E35', '  +0.56%  ')
    ,@('D36', '1.02')
    ,@('E36', '  +17.65%  ')
    ,@('D37', '4.06')
    ,@('E37', '  +4.72%  ')
    ,@('E38', '  +4.46%  ')
    ,@('D39', '36.89')
    ,@('E39', '  +3.29%  ')
    ,@('D40', '0.850')
    ,@('E40', '  +3.90%  ')
    ,@('D41', '3.70')
    ,@('E41', '  +6.50%  ')
    ,@('D42', '1.42')
    ,@('E42', '  +3.01%  ')
    ,@('D43', '0.623')
    ,@('E43', '  +2.95%  ')
    ,@('D44', '281.58')
    ,@('E44', '  -1.01%  ')
    ,@('B45', 'EnergySwap')
    ,@('C45', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens')
    ,@('D45', '19.98')
    ,@('E45', '  +7.52%  ')
    ,@('B46', 'FirstDigitalUSD')
    ,@('C46', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd')
    ,@('D46', '0.992')
    ,@('E46', '  -0.69%  ')
    ,@('B47', 'Stellar')
    ,@('C47', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm')
    ,@('D47', '0.0983')
    ,@('E47', '  +0.47%  ')
    ,@('D48', '0.0535')
    ,@('E48', '  +1.53%  ')
    ,@('B49', 'RenderToken')
    ,@('C49', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr')
    ,@('D49', '4.76')
    ,@('E49', '  +5.90%  ')
    ,@('B50', 'Maker')
    ,@('C50', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr')
    ,@('D50', '2.014.88')
    ,@('E50', '  +5.71%  ')
    ,@('B51', 'VeChain')
    ,@('C51', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet')
    ,@('D51', '0.0231')
    ,@('E51', '  +2.72%  ')
)

foreach ($u in $updates) {
    Set-TextValue $u[0] $u[1]
}

Write-Host "Applied $($updates.Count) cell updates"